$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

$row = 12

# Column A holds a numeric-looking issue id that must stay text (like the
# existing rows), so force text formatting before assigning the value.
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("A" + $row).Value = "16"

$ws.Range("B" + $row).Value = "newrelic_alert_condition and newrelic_cloud_aws_integrations"
$ws.Range("C" + $row).Value = "open"
$ws.Range("D" + $row).Value = "2025-03-24T10:26:30Z"
$ws.Range("E" + $row).Value = "bug"
